$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.45279999999999
$ws.Range("A9").Value = -21.43840000000001
$ws.Range("C9").Value = -10.2817
$ws.Range("C11").Value = -12.91639999999999
$ws.Range("A13").Value = -22.2827
$ws.Range("A16").Value = -21.70239999999999
$ws.Range("C16").Value = -12.8779
$ws.Range("A18").Value = -22.24830000000001
$ws.Range("A20").Value = -21.50829999999998
$ws.Range("C23").Value = -11.9959
$ws.Range("C24").Value = -12.4064
$ws.Range("A26").Value = -21.92499999999998
$ws.Range("C26").Value = -12.89109999999999
$ws.Range("A27").Value = -21.88739999999998
$ws.Range("A29").Value = -21.54879999999998
$ws.Range("C34").Value = -11.91400000000001
$ws.Range("A35").Value = -21.80319999999998
$ws.Range("C35").Value = -12.5682
$ws.Range("A36").Value = -20.67089999999998
$ws.Range("C44").Value = -12.84939999999999
$ws.Range("A45").Value = -21.46009999999999
$ws.Range("C48").Value = -12.1701
$ws.Range("C49").Value = -14.48839999999999
$ws.Range("C52").Value = -11.153
$ws.Range("A55").Value = -22.20850000000001
$ws.Range("A57").Value = -22.31160000000001
$ws.Range("C66").Value = -11.1812
$ws.Range("C67").Value = -11.22720000000001
$ws.Range("A69").Value = -21.54029999999998
$ws.Range("C73").Value = -11.26390000000001
$ws.Range("A76").Value = -19.68169999999999
$ws.Range("A78").Value = -20.5596
$ws.Range("C78").Value = -11.7799
$ws.Range("C80").Value = -11.8522
$ws.Range("A82").Value = -21.84840000000001
$ws.Range("A83").Value = -21.58249999999999
$ws.Range("C91").Value = -12.3738
$ws.Range("A93").Value = -21.43300000000002
$ws.Range("A97").Value = -21.67970000000001
$ws.Range("C97").Value = -11.04900000000001
$ws.Range("C99").Value = -12.696
$ws.Range("C104").Value = -12.90090000000001
